$wb = $excel.ActiveWorkbook

# Sheet1 (展览): update F-column 'want to go' counters
$wsExpo = $wb.Worksheets.Item(1)
$wsExpo.Range("F3").Value = 1734
$wsExpo.Range("F4").Value = 97
$wsExpo.Range("F5").Value = 62
$wsExpo.Range("F6").Value = 698
$wsExpo.Range("F7").Value = 172
$wsExpo.Range("F8").Value = 197
$wsExpo.Range("F11").Value = 30
$wsExpo.Range("F12").Value = 560
$wsExpo.Range("F13").Value = 482
$wsExpo.Range("F16").Value = 133
$wsExpo.Range("F17").Value = 774
$wsExpo.Range("F18").Value = 2577
$wsExpo.Range("F23").Value = 186
$wsExpo.Range("F25").Value = 130
$wsExpo.Range("F27").Value = 941
$wsExpo.Range("F29").Value = 154
$wsExpo.Range("F33").Value = 254

# Sheet3 (本地生活): update F-column 'want to go' counters
$wsLocal = $wb.Worksheets.Item(3)
$wsLocal.Range("F2").Value = 1760
$wsLocal.Range("F4").Value = 41
$wsLocal.Range("F5").Value = 2347
$wsLocal.Range("F6").Value = 948
$wsLocal.Range("F9").Value = 1190
$wsLocal.Range("F10").Value = 304

# Sheet4 (全部类型): update F-column 'want to go' counters
$wsAll = $wb.Worksheets.Item(4)
$wsAll.Range("F2").Value = 1760
$wsAll.Range("F4").Value = 41
$wsAll.Range("F5").Value = 2347
$wsAll.Range("F7").Value = 1734
$wsAll.Range("F9").Value = 948
$wsAll.Range("F10").Value = 1190
$wsAll.Range("F11").Value = 304
$wsAll.Range("F13").Value = 97
$wsAll.Range("F14").Value = 62
$wsAll.Range("F15").Value = 698
$wsAll.Range("F16").Value = 172
$wsAll.Range("F18").Value = 197
$wsAll.Range("F20").Value = 30
$wsAll.Range("F21").Value = 560
$wsAll.Range("F22").Value = 482
$wsAll.Range("F25").Value = 133
$wsAll.Range("F26").Value = 774
$wsAll.Range("F27").Value = 2577
$wsAll.Range("F32").Value = 186
$wsAll.Range("F33").Value = 130
$wsAll.Range("F35").Value = 941
$wsAll.Range("F36").Value = 524
$wsAll.Range("F37").Value = 81
$wsAll.Range("F39").Value = 154
$wsAll.Range("F44").Value = 260
$wsAll.Range("F49").Value = 254

# Sheet2 (演出): remove first data row (2024-02-04 event) by shifting rows 3-32
# up by one (values only, column A/index stays as-is), then delete now-empty last row 33
$wsPerf = $wb.Worksheets.Item(2)

# Row 3
$wsPerf.Range("B3").Value = '2024-02-14'
$wsPerf.Range("C3").Value = '上海·【情人节特辑】《那年我们》记忆重启韩剧经典OST音乐会《请回答1988》《来自星星的你》（取消）'
$wsPerf.Range("D3").Value = '牛庄路704号 中国大戏院'
$wsPerf.Range("E3").Value = '2024.02.14 19:30-02.14 21:00'
$wsPerf.Range("F3").Value = 7
$wsPerf.Range("G3").Value = '不可售'
$wsPerf.Range("H3").Value = 'https://show.bilibili.com/platform/detail.html?id=80615'
$wsPerf.Range("I3").Value = '//i0.hdslb.com/bfs/openplatform/202401/5DDVhKcO1704767761361.png'

# Row 4
$wsPerf.Range("B4").Value = '2024-02-20'
$wsPerf.Range("C4").Value = '上海·Liyuu 「鲤好！」粉丝见面会'
$wsPerf.Range("D4").Value = '宜昌路179号 万代南梦宫上海文化中心'
$wsPerf.Range("E4").Value = '2024.02.20 14:30-02.20 21:00'
$wsPerf.Range("F4").Value = 286
$wsPerf.Range("G4").Value = 360
$wsPerf.Range("H4").Value = 'https://show.bilibili.com/platform/detail.html?id=81740'
$wsPerf.Range("I4").Value = '//i1.hdslb.com/bfs/openplatform/202402/LywKSi4B1707040250585.png'

# Row 5
$wsPerf.Range("B5").Value = '2024-02-23'
$wsPerf.Range("C5").Value = '上海·天空之城-经典动漫烛光音乐会'
$wsPerf.Range("D5").Value = '曹杨路1888号 上海露边社·演艺空间'
$wsPerf.Range("E5").Value = '2024.02.23 19:30-02.23 21:00'
$wsPerf.Range("F5").Value = 2
$wsPerf.Range("G5").Value = 88
$wsPerf.Range("H5").Value = 'https://show.bilibili.com/platform/detail.html?id=81541'
$wsPerf.Range("I5").Value = '//i2.hdslb.com/bfs/openplatform/202402/Q3L80ixO1706778157039.jpeg'

# Row 6
$wsPerf.Range("B6").Value = '2024-02-24'
$wsPerf.Range("C6").Value = '上海·《哈利的魔法世界》动漫视听音乐会'
$wsPerf.Range("D6").Value = '都市路4889号（莘庄地铁站南广场） 上海保利城市剧院'
$wsPerf.Range("E6").Value = '2024.02.24 14:30-02.24 16:00'
$wsPerf.Range("F6").Value = 13
$wsPerf.Range("G6").Value = 158
$wsPerf.Range("H6").Value = 'https://show.bilibili.com/platform/detail.html?id=80639'
$wsPerf.Range("I6").Value = '//i2.hdslb.com/bfs/openplatform/202401/4PieCC9N1706261750579.jpeg'

# Row 7
$wsPerf.Range("B7").Value = '2024-02-25'
$wsPerf.Range("C7").Value = '上海·青山吉能见面会'
$wsPerf.Range("D7").Value = '虹许路731号4号楼 THE BOXX•城市乐园'
$wsPerf.Range("E7").Value = '2024.02.25 14:30-02.25 19:30'
$wsPerf.Range("F7").Value = 219
$wsPerf.Range("G7").Value = 380
$wsPerf.Range("H7").Value = 'https://show.bilibili.com/platform/detail.html?id=80142'
$wsPerf.Range("I7").Value = '//i0.hdslb.com/bfs/openplatform/202312/1npuHFBM1703231674558.jpeg'

# Row 8
$wsPerf.Range("B8").Value = '2024-03-02'
$wsPerf.Range("C8").Value = '上海·2024藤田玲上海粉丝见面会'
$wsPerf.Range("D8").Value = '宜昌路179号 万代南梦宫上海文化中心'
$wsPerf.Range("E8").Value = '2024.03.02 12:30-03.02 19:40'
$wsPerf.Range("F8").Value = 15
$wsPerf.Range("G8").Value = 580
$wsPerf.Range("H8").Value = 'https://show.bilibili.com/platform/detail.html?id=80993'
$wsPerf.Range("I8").Value = '//i0.hdslb.com/bfs/openplatform/202401/Vm6ntgVd1705548188785.png'

# Row 9
$wsPerf.Range("B9").Value = '2024-03-02'
$wsPerf.Range("C9").Value = '上海·小山百代2024上海粉丝见面会'
$wsPerf.Range("D9").Value = '宜昌路179号 万代南梦宫上海文化中心'
$wsPerf.Range("E9").Value = '2024.03.02 13:00-03.02 20:00'
$wsPerf.Range("F9").Value = 292
$wsPerf.Range("G9").Value = 380
$wsPerf.Range("H9").Value = 'https://show.bilibili.com/platform/detail.html?id=80924'
$wsPerf.Range("I9").Value = '//i1.hdslb.com/bfs/openplatform/202401/FpA9OkKy1705467080070.jpeg'

# Row 10
$wsPerf.Range("B10").Value = '2024-03-08'
$wsPerf.Range("C10").Value = '上海·《月亮代表我的心》摇滚情歌之夜--630乐团演绎经典'
$wsPerf.Range("D10").Value = '淞沪路388号创智天地广场7号楼一层 创智天地梦剧场'
$wsPerf.Range("E10").Value = '2024.03.08 20:00-03.08 21:30'
$wsPerf.Range("F10").Value = 0
$wsPerf.Range("G10").Value = 90
$wsPerf.Range("H10").Value = 'https://show.bilibili.com/platform/detail.html?id=81676'
$wsPerf.Range("I10").Value = '//i1.hdslb.com/bfs/openplatform/202402/0yMuaTCo1706860617422.png'

# Row 11
$wsPerf.Range("B11").Value = '2024-03-09'
$wsPerf.Range("C11").Value = '上海·《挪威的森林》—摇滚情歌之夜演唱会'
$wsPerf.Range("D11").Value = '南京西路1376号 上海商城剧院'
$wsPerf.Range("E11").Value = '2024.03.09 19:30-03.09 21:00'
$wsPerf.Range("F11").Value = 0
$wsPerf.Range("G11").Value = 72
$wsPerf.Range("H11").Value = 'https://show.bilibili.com/platform/detail.html?id=81241'
$wsPerf.Range("I11").Value = '//i2.hdslb.com/bfs/openplatform/202401/1FJ0Fj5m1705915336335.jpeg'

# Row 12
$wsPerf.Range("B12").Value = '2024-03-09'
$wsPerf.Range("C12").Value = '上海·爱乐之城音乐会'
$wsPerf.Range("D12").Value = '南京西路1376号 上海商城剧院'
$wsPerf.Range("E12").Value = '2024.03.09 14:00-03.09 15:30'
$wsPerf.Range("F12").Value = 5
$wsPerf.Range("G12").Value = 60
$wsPerf.Range("H12").Value = 'https://show.bilibili.com/platform/detail.html?id=81289'
$wsPerf.Range("I12").Value = '//i2.hdslb.com/bfs/openplatform/202401/ZZXtDrwZ1705996679699.jpeg'

# Row 13
$wsPerf.Range("B13").Value = '2024-03-10'
$wsPerf.Range("C13").Value = '上海·三森铃子10周年纪念2024演唱会'
$wsPerf.Range("D13").Value = '宜昌路179号 万代南梦宫上海文化中心'
$wsPerf.Range("E13").Value = '2024.03.10 18:00-03.10 19:30'
$wsPerf.Range("F13").Value = 524
$wsPerf.Range("G13").Value = 399
$wsPerf.Range("H13").Value = 'https://show.bilibili.com/platform/detail.html?id=81433'
$wsPerf.Range("I13").Value = '//i0.hdslb.com/bfs/openplatform/202401/L8rmm2h81706236781799.jpeg'

# Row 14
$wsPerf.Range("B14").Value = '2024-03-16'
$wsPerf.Range("C14").Value = '上海·三月的幻想演唱会2024「飞越蓝色时刻」'
$wsPerf.Range("D14").Value = '宜昌路179号 万代南梦宫上海文化中心'
$wsPerf.Range("E14").Value = '2024.03.16 19:00-03.16 20:30'
$wsPerf.Range("F14").Value = 81
$wsPerf.Range("G14").Value = 380
$wsPerf.Range("H14").Value = 'https://show.bilibili.com/platform/detail.html?id=80811'
$wsPerf.Range("I14").Value = '//i1.hdslb.com/bfs/openplatform/202401/TO6xpSqr1705289483473.png'

# Row 15
$wsPerf.Range("B15").Value = '2024-03-17'
$wsPerf.Range("C15").Value = '上海 ·《疯狂动物城》动漫视听音乐会'
$wsPerf.Range("D15").Value = '牛庄路704号 中国大戏院'
$wsPerf.Range("E15").Value = '2024.03.17 15:30-03.17 17:00'
$wsPerf.Range("F15").Value = 9
$wsPerf.Range("G15").Value = 80
$wsPerf.Range("H15").Value = 'https://show.bilibili.com/platform/detail.html?id=81112'
$wsPerf.Range("I15").Value = '//i2.hdslb.com/bfs/openplatform/202401/Wg8b6SRn1705651166088.png'

# Row 16
$wsPerf.Range("B16").Value = '2024-03-17'
$wsPerf.Range("C16").Value = '上海·amazarashi Asia Tour 2024 「永遠市 -Eternal City-」上海公演'
$wsPerf.Range("D16").Value = '宜昌路179号 万代南梦宫上海文化中心'
$wsPerf.Range("E16").Value = '2024.03.17 18:00-03.17 19:30'
$wsPerf.Range("F16").Value = 951
$wsPerf.Range("G16").Value = '已售罄'
$wsPerf.Range("H16").Value = 'https://show.bilibili.com/platform/detail.html?id=81039'
$wsPerf.Range("I16").Value = '//i2.hdslb.com/bfs/openplatform/202401/icsawZU11705566039011.jpeg'

# Row 17
$wsPerf.Range("B17").Value = '2024-03-17'
$wsPerf.Range("C17").Value = '上海·《笑傲江湖》经典武侠影视金曲音乐会'
$wsPerf.Range("D17").Value = '牛庄路704号 中国大戏院'
$wsPerf.Range("E17").Value = '2024.03.17 19:30-03.17 21:00'
$wsPerf.Range("F17").Value = 1
$wsPerf.Range("G17").Value = 80
$wsPerf.Range("H17").Value = 'https://show.bilibili.com/platform/detail.html?id=80875'
$wsPerf.Range("I17").Value = '//i1.hdslb.com/bfs/openplatform/202401/8AwIAy4I1705385447242.jpeg'

# Row 18
$wsPerf.Range("B18").Value = '2024-03-17'
$wsPerf.Range("C18").Value = '上海·遇见新海诚--帝玖「这次一定」室内乐ACG音乐会'
$wsPerf.Range("D18").Value = '延安东路523号 凯迪拉克·上海音乐厅'
$wsPerf.Range("E18").Value = '2024.03.17 14:00-03.17 16:00'
$wsPerf.Range("F18").Value = 30
$wsPerf.Range("G18").Value = 80
$wsPerf.Range("H18").Value = 'https://show.bilibili.com/platform/detail.html?id=81258'
$wsPerf.Range("I18").Value = '//i1.hdslb.com/bfs/openplatform/202401/eysvN81k1705977896972.jpeg'

# Row 19
$wsPerf.Range("B19").Value = '2024-03-21'
$wsPerf.Range("C19").Value = '上海·春卷饭 十周年  2024  专场演出'
$wsPerf.Range("D19").Value = '嘉兴路街道瑞虹路188号瑞虹天地月亮湾3层 Modern Sky LAB摩登天空(瑞虹天地店)'
$wsPerf.Range("E19").Value = '2024.03.21 20:00-03.21 22:00'
$wsPerf.Range("F19").Value = 605
$wsPerf.Range("G19").Value = '已售罄'
$wsPerf.Range("H19").Value = 'https://show.bilibili.com/platform/detail.html?id=81190'
$wsPerf.Range("I19").Value = '//i1.hdslb.com/bfs/openplatform/202401/ho9rIMg21705894649801.jpeg'

# Row 20
$wsPerf.Range("B20").Value = '2024-03-23'
$wsPerf.Range("C20").Value = '上海·《卡农Canon in D》世界经典作品视听音乐会'
$wsPerf.Range("D20").Value = '南京西路1376号 上海商城剧院'
$wsPerf.Range("E20").Value = '2024.03.23 19:30-03.23 21:00'
$wsPerf.Range("F20").Value = 1
$wsPerf.Range("G20").Value = 50
$wsPerf.Range("H20").Value = 'https://show.bilibili.com/platform/detail.html?id=81358'
$wsPerf.Range("I20").Value = '//i1.hdslb.com/bfs/openplatform/202401/Ctne29Xn1706089385959.png'

# Row 21
$wsPerf.Range("B21").Value = '2024-03-23'
$wsPerf.Range("C21").Value = '上海·《四月是你的谎言》友人A经典动漫音乐会'
$wsPerf.Range("D21").Value = '南京西路1376号 上海商城剧院'
$wsPerf.Range("E21").Value = '2024.03.23 15:00-03.23 16:30'
$wsPerf.Range("F21").Value = 31
$wsPerf.Range("G21").Value = 50
$wsPerf.Range("H21").Value = 'https://show.bilibili.com/platform/detail.html?id=81361'
$wsPerf.Range("I21").Value = '//i0.hdslb.com/bfs/openplatform/202401/wL0ZWVYi1706091574963.png'

# Row 22
$wsPerf.Range("B22").Value = '2024-03-24'
$wsPerf.Range("C22").Value = '上海·“燃魂巅峰交响版”VICTORY·星球大战·加勒比海盗 大型交响音乐会'
$wsPerf.Range("D22").Value = '丁香路425号(上海科技馆地铁站1号口步行460米) 上海东方艺术中心音乐厅'
$wsPerf.Range("E22").Value = '2024.03.24 19:30-03.24 21:00'
$wsPerf.Range("F22").Value = 1
$wsPerf.Range("G22").Value = 80
$wsPerf.Range("H22").Value = 'https://show.bilibili.com/platform/detail.html?id=81501'
$wsPerf.Range("I22").Value = '//i2.hdslb.com/bfs/openplatform/202401/IEM4vSmT1706520953088.jpeg'

# Row 23
$wsPerf.Range("B23").Value = '2024-03-29'
$wsPerf.Range("C23").Value = '上海·KANAKO ITO&AYANE 2024 LIVE'
$wsPerf.Range("D23").Value = '宜昌路179号 万代南梦宫上海文化中心'
$wsPerf.Range("E23").Value = '2024.03.29 19:00-03.29 20:30'
$wsPerf.Range("F23").Value = 260
$wsPerf.Range("G23").Value = 380
$wsPerf.Range("H23").Value = 'https://show.bilibili.com/platform/detail.html?id=81416'
$wsPerf.Range("I23").Value = '//i0.hdslb.com/bfs/openplatform/202401/4Y4U8tC01706172039039.jpeg'

# Row 24
$wsPerf.Range("B24").Value = '2024-03-30'
$wsPerf.Range("C24").Value = '上海· TRUE（唐沢美帆）上海动漫交响音乐会'
$wsPerf.Range("D24").Value = '丁香路425号 上海东方艺术中心'
$wsPerf.Range("E24").Value = '2024.03.30 19:30-03.30 21:00'
$wsPerf.Range("F24").Value = 235
$wsPerf.Range("G24").Value = 680
$wsPerf.Range("H24").Value = 'https://show.bilibili.com/platform/detail.html?id=80906'
$wsPerf.Range("I24").Value = '//i0.hdslb.com/bfs/openplatform/202401/FaJbLvS51705401178235.jpeg'

# Row 25
$wsPerf.Range("B25").Value = '2024-03-31'
$wsPerf.Range("C25").Value = '上海·《天空之城》宫崎骏&久石让经典作品动漫视听音乐会'
$wsPerf.Range("D25").Value = '江宁路466号 上海艺海剧院·小剧场'
$wsPerf.Range("E25").Value = '2024.03.31 10:30-03.31 12:00'
$wsPerf.Range("F25").Value = 1
$wsPerf.Range("G25").Value = 50
$wsPerf.Range("H25").Value = 'https://show.bilibili.com/platform/detail.html?id=81660'
$wsPerf.Range("I25").Value = '//i2.hdslb.com/bfs/openplatform/202402/QKmfdsEM1706853934802.jpeg'

# Row 26
$wsPerf.Range("B26").Value = '2024-03-31'
$wsPerf.Range("C26").Value = '上海·《热血之巅·突破次元壁》ACG动漫电影音乐会'
$wsPerf.Range("D26").Value = '江宁路466号 上海艺海剧院·小剧场'
$wsPerf.Range("E26").Value = '2024.03.31 15:00-03.31 16:30'
$wsPerf.Range("F26").Value = 2
$wsPerf.Range("G26").Value = 90
$wsPerf.Range("H26").Value = 'https://show.bilibili.com/platform/detail.html?id=81672'
$wsPerf.Range("I26").Value = '//i2.hdslb.com/bfs/openplatform/202402/5k9iIwRO1706859635834.jpeg'

# Row 27
$wsPerf.Range("B27").Value = '2024-04-06'
$wsPerf.Range("C27").Value = '上海·从Butter-Fly到夏目之爱してる —— “好想大声说爱你”动漫钢琴演奏会'
$wsPerf.Range("D27").Value = '复兴中路1380号 捷豹上海交响音乐厅'
$wsPerf.Range("E27").Value = '2024.04.06 19:30-04.06 21:30'
$wsPerf.Range("F27").Value = 14
$wsPerf.Range("G27").Value = 80
$wsPerf.Range("H27").Value = 'https://show.bilibili.com/platform/detail.html?id=80050'
$wsPerf.Range("I27").Value = '//i0.hdslb.com/bfs/openplatform/202312/0iJP3TY61703056498448.jpeg'

# Row 28
$wsPerf.Range("B28").Value = '2024-04-13'
$wsPerf.Range("C28").Value = '上海·《四月是你的谎言》——“公生”与“薰”的钢琴小提琴唯美经典音乐集'
$wsPerf.Range("D28").Value = '丁香路425号 上海东方艺术中心'
$wsPerf.Range("E28").Value = '2024.04.13 19:30-04.13 21:30'
$wsPerf.Range("F28").Value = 190
$wsPerf.Range("G28").Value = 80
$wsPerf.Range("H28").Value = 'https://show.bilibili.com/platform/detail.html?id=78667'
$wsPerf.Range("I28").Value = '//i1.hdslb.com/bfs/openplatform/202311/bTP7w6GD1700130122940.jpeg'

# Row 29
$wsPerf.Range("B29").Value = '2024-04-20'
$wsPerf.Range("C29").Value = '上海·Laurent Coulondre“心动巴黎”2024中国巡回音乐会'
$wsPerf.Range("D29").Value = '汾阳路20号上海音乐学院内 上海贺绿汀音乐厅'
$wsPerf.Range("E29").Value = '2024.04.20 19:30-04.20 21:30'
$wsPerf.Range("F29").Value = 4
$wsPerf.Range("G29").Value = 80
$wsPerf.Range("H29").Value = 'https://show.bilibili.com/platform/detail.html?id=81135'
$wsPerf.Range("I29").Value = '//i2.hdslb.com/bfs/openplatform/202401/wXDdS5ap1705651730828.jpeg'

# Row 30
$wsPerf.Range("B30").Value = '2024-04-26'
$wsPerf.Range("C30").Value = '上海· 夏川里美 2024 巡回演唱会 出道 25 周年纪念专场'
$wsPerf.Range("D30").Value = '东大名路889号 友邦大剧院'
$wsPerf.Range("E30").Value = '2024.04.26 19:30-04.26 21:30'
$wsPerf.Range("F30").Value = 23
$wsPerf.Range("G30").Value = 280
$wsPerf.Range("H30").Value = 'https://show.bilibili.com/platform/detail.html?id=81139'
$wsPerf.Range("I30").Value = '//i2.hdslb.com/bfs/openplatform/202401/0Fj4cYOH1705652393930.jpeg'

# Row 31
$wsPerf.Range("B31").Value = '2024-05-17'
$wsPerf.Range("C31").Value = '上海·Rie fu 船越里惠 日本知名唱作歌手2024出道20周年中国巡回演唱会'
$wsPerf.Range("D31").Value = '南京西路1376号 上海商城剧院'
$wsPerf.Range("E31").Value = '2024.05.17 19:30-05.17 21:00'
$wsPerf.Range("F31").Value = 2
$wsPerf.Range("G31").Value = 380
$wsPerf.Range("H31").Value = 'https://show.bilibili.com/platform/detail.html?id=81506'
$wsPerf.Range("I31").Value = '//i2.hdslb.com/bfs/openplatform/202401/6ue4xoaR1706523724335.jpeg'

# Row 32
$wsPerf.Range("B32").Value = '2024-06-08'
$wsPerf.Range("C32").Value = '上海·菊次郎的夏天——久石让钢琴曲梦幻之旅演奏会'
$wsPerf.Range("D32").Value = '延安东路523号 凯迪拉克·上海音乐厅'
$wsPerf.Range("E32").Value = '2024.06.08 19:30-06.08 21:00'
$wsPerf.Range("F32").Value = 6
$wsPerf.Range("G32").Value = 80
$wsPerf.Range("H32").Value = 'https://show.bilibili.com/platform/detail.html?id=81413'
$wsPerf.Range("I32").Value = '//i2.hdslb.com/bfs/openplatform/202401/QqKuy4611706169245363.jpeg'

# Remove the now-duplicated last row (33) entirely
$wsPerf.Rows.Item(33).Delete()
